# Apply the "Added additional items in ToInitialize.xlsx" edit.
# This adds new Filename/StrFind/StrReplace rows (new source files whose
# generated C code needs extra find/replace initialization entries),
# drops the trailing semicolons from the tclust_wrapper1.c detpar/rotpar/
# shapepar StrReplace values, and pushes the trailing printf/fflush rows
# further down the table, followed by two new HAwei.c rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 21-23: tclust_wrapper1.c block loses the trailing ';' on the
#     StrReplace (column C) values ------------------------------------
$ws.Range("C21").Value = "double detpar=0"
$ws.Range("C22").Value = "double rotpar=0"
$ws.Range("C23").Value = "double shapepar=0"

# --- Rows 24-41: brand new Filename/StrFind/StrReplace triples --------
$newRows = @(
    @(24, "HArho.c",        "double c_data",    "double c_data=0"),
    @(25, "HArho.c",        "double b_data",    "double b_data=0"),
    @(26, "HArho.c",        "double a_data",    "double a_data=0"),
    @(27, "MMregcore.c",    "double c2;",        "double c2=0;"),
    @(28, "MMregcore.c",    "double b2;",        "double b2=0;"),
    @(29, "MMregcore.c",    "double a2",         "double a2=0"),
    @(30, "MMregcore.c",    "double A2",         "double A2=0"),
    @(31, "Mscale.c",       "double scnew",      "double scnew=0"),
    @(32, "OPTwei.c",       "double x1_data",    "double x1_data=0"),
    @(33, "Sreg.c",         "double scaletest",  "double scaletest=0"),
    @(34, "Sreg.c",         "double kc",         "double kc=0"),
    @(35, "Sreg.c",         "double A",          "double A=0"),
    @(36, "Sreg.c",         "double c;",         "double c=0;"),
    @(37, "Sreg.c",         "double d",          "double d=0"),
    @(38, "Sreg_wrapper.c", "double scaletest",  "double scaletest=0"),
    @(39, "Sreg_wrapper.c", "double kc",         "double kc=0"),
    @(40, "Sreg_wrapper.c", "double A",          "double A=0"),
    @(41, "Sreg_wrapper.c", "double B",          "double B=0")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
}

# --- Rows 42-43: the printf(/Rprintf( and fflush(stdout/ //fflush(stdout
#     rows shift two rows further down (no column A value) -------------
$ws.Range("B42").Value = "printf("
$ws.Range("C42").Value = "Rprintf("
$ws.Range("B43").Value = "fflush(stdout"
$ws.Range("C43").Value = "//fflush(stdout"

# --- Rows 44-45: two new HAwei.c rows appended at the end -------------
$ws.Range("A44").Value = "HAwei.c"
$ws.Range("B44").Value = "double c_data;"
$ws.Range("C44").Value = "double c_data=0;"

$ws.Range("A45").Value = "HAwei.c"
$ws.Range("B45").Value = "double b_data;"
$ws.Range("C45").Value = "double b_data=0;"

# --- Update the sheet selection/view to match the new active cell -----
$ws.Range("B46").Select()
